# Updated cryptos list on Sat Dec 30 10:08:52 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.917.54"
$ws.Range("E2").Value = "  -2.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.281.59"
$ws.Range("E3").Value = "  -3.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "315.32"
$ws.Range("E5").Value = "  -0.72%  "

# Row 6 - Solana
$ws.Range("D6").Value = "102.46"
$ws.Range("E6").Value = "  -6.07%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.55%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -3.19%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "38.51"
$ws.Range("E10").Value = "  -8.03%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0902"
$ws.Range("E11").Value = "  -2.66%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "8.19"
$ws.Range("E12").Value = "  -4.80%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.54%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.953"
$ws.Range("E14").Value = "  -5.11%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "15.21"
$ws.Range("E15").Value = "  -5.75%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.625.74"
$ws.Range("E16").Value = "  -3.62%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.279.08"
$ws.Range("E17").Value = "  -4.22%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "41.812.01"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "7.43"
$ws.Range("E19").Value = "  -3.93%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -1.66%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "73.12"
$ws.Range("E21").Value = "  -4.24%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "277.31"
$ws.Range("E22").Value = "  +8.03%  "

# Row 23 - PancakeSwap
$ws.Range("D23").Value = "3.54"
$ws.Range("E23").Value = "  -3.05%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "10.04"
$ws.Range("E24").Value = "  +5.40%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -3.33%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.74%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "10.68"
$ws.Range("E27").Value = "  -7.13%  "

# Row 28 - Toncoin
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +3.44%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "22.86"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30 - Monero
$ws.Range("D30").Value = "162.67"
$ws.Range("E30").Value = "  -5.26%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "34.60"
$ws.Range("E31").Value = "  -7.61%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "0.0865"
$ws.Range("E32").Value = "  -3.19%  "

# Row 33 - WEMIXToken
$ws.Range("D33").Value = "2.87"
$ws.Range("E33").Value = "  -0.97%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  -4.47%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  -0.01%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -7.17%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").Value = "  -3.52%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  +6.64%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.0345"
$ws.Range("E39").Value = "  -5.42%  "

# Row 40 - NEARProtocol
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -7.62%  "

# Row 41 - BitcoinSV
$ws.Range("D41").Value = "100.08"
$ws.Range("E41").Value = "  +16.83%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  -3.65%  "

# Row 43 - MultiversX
$ws.Range("D43").Value = "68.85"
$ws.Range("E43").Value = "  -3.27%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  +0.09%  "

# Row 45 - Algorand
$ws.Range("E45").Value = "  -8.00%  "

# Row 46 - Aave
$ws.Range("D46").Value = "115.57"
$ws.Range("E46").Value = "  +2.69%  "

# Row 47 - Celestia
$ws.Range("D47").Value = "11.77"
$ws.Range("E47").Value = "  -4.41%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "8.95"
$ws.Range("E48").Value = "  -3.15%  "

# Row 49 and 50 swap places: THORChain <-> ordi
# Row 49 becomes ordi
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "75.39"
$ws.Range("E49").Value = "  -2.50%  "

# Row 50 becomes THORChain
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  -6.18%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  -3.95%  "
